$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("#system")
# try applying existing used style cell's format via Copy  
$ws2 = $wb.Worksheets.Item("macros")
Write-Host $ws2.Range("A1").Value2
